# timesheet.xlsx edit: add new columns (supervisor/cost centre/school/start
# date/"added") to Sheet1, replacing the old "WEEKLY TOTAL" column, and wire
# up a fresh sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column N used to hold "WEEKLY TOTAL " as a header with no data underneath
# (it was effectively unused). Rename it and clear the data column, then add
# four brand-new header-only columns after it.
$ws.Cells.Item(1, 14).Value = "SUPERVIOSR NAME"
$ws.Cells.Item(1, 15).Value = "COST CENTRE"
$ws.Cells.Item(1, 16).Value = "SCHOOL NAME"
$ws.Cells.Item(1, 17).Value = "START DATE"
$ws.Cells.Item(1, 18).Value = "ADDED"

$ws.Range("M1").Copy()
$ws.Range("N1:R1").PasteSpecial(-4122)  # xlPasteFormats

# Move/resize the view: scroll so column B is left-most and select the new
# "ADDED" header cell.
$ws.Range("R2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2

$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(12).AutoFit()
$ws.Columns.Item(14).AutoFit()
$ws.Columns.Item(15).AutoFit()
$ws.Columns.Item(16).AutoFit()
$ws.Columns.Item(17).AutoFit()
